# Fix xpath expressions: the "reviews_count" column (E) was removed from
# the data source, so delete the corresponding column in the header row.
# Deleting the entire column shifts all subsequent columns (F:K) one
# position to the left and updates the sheet dimension accordingly.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Columns.Item(5).Delete()
